$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.848.28"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "1.649.37"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "'217.12"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "'0.503"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'0.0629"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'19.28"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "'0.0845"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "1.638.40"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "'4.18"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "'0.530"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "'64.92"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "26.807.47"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "0.0₃0738"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "'215.22"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").Value = "'4.40"
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("E21").Value = "  +12.66%  "
$ws.Range("D22").Value = "'6.28"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "'9.39"
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("D24").Value = "'147.21"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").Value = "'7.21"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").Value = "'15.71"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "'0.0511"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").Value = "'3.37"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").Value = "'3.02"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").Value = "1.296.29"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("D36").Value = "'0.0177"
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("D37").Value = "'0.538"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").Value = "'0.825"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").Value = "'0.808"
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").Value = "'5.33"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").Value = "1.783.63"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "'62.13"
$ws.Range("E44").Value = "  +3.69%  "
$ws.Range("D45").Value = "'91.96"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'1.62"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("D48").Value = "'0.0522"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").Value = "'7.65"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").Value = "'0.0975"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").Value = "'0.408"
$ws.Range("E51").Value = "  +0.39%  "
